# Apply edit: correct the "Latest period (release date)" value for the
# "Employment share by occupation" row (C3), and update the active
# selection to C4, matching the author's recorded change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the data value in C3 (was "Jul 2021 - Jun 2022 (11/10/22)")
$ws.Range("C3").Value = "Jan 2021 - Dec 2021 (12/04/22)"

# Update the saved selection to match the author's cursor position
$ws.Range("C4").Select()
